$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = 24  # data rows 2..25

# --- Update existing columns I, K, M, O (swap 1<->2 values) for rows 2-25 ---
$iko = New-Object 'object[,]' $rows,7
for ($i = 0; $i -lt $rows; $i++) {
    $iko[$i,0] = 2  # I: was 1
    $iko[$i,1] = 2  # J: unchanged
    $iko[$i,2] = 1  # K: was 2
    $iko[$i,3] = 2  # L: unchanged
    $iko[$i,4] = 2  # M: was 1
    $iko[$i,5] = 2  # N: unchanged
    $iko[$i,6] = 1  # O: was 2
}
$ws.Range("I2:O25").Value = $iko

# --- New columns P and Q for data rows 2-25, value 2 ---
$pq = New-Object 'object[,]' $rows,2
for ($i = 0; $i -lt $rows; $i++) {
    $pq[$i,0] = 2
    $pq[$i,1] = 2
}
$ws.Range("P2:Q25").Value = $pq

# --- New header cells P1 = 14, Q1 = 15, matching style of existing header cells ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
